$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '93.827.19'
$ws.Range('E2').Value = '  -4.12%  '
$ws.Range('D3').Value = '3.404.53'
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '236.40'
$ws.Range('E5').Value = '  -6.82%  '
$ws.Range('D6').Value = '638.08'
$ws.Range('E6').Value = '  -3.59%  '
$ws.Range('D7').Value = '1.42'
$ws.Range('E7').Value = '  -4.09%  '
$ws.Range('D8').Value = '0.399'
$ws.Range('E8').Value = '  -6.66%  '
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('D10').Value = '0.960'
$ws.Range('E10').Value = '  -7.90%  '
$ws.Range('D11').Value = '3.402.81'
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('D12').Value = '0.198'
$ws.Range('E12').Value = '  -5.06%  '
$ws.Range('D13').Value = '41.21'
$ws.Range('E13').Value = '  -7.13%  '
$ws.Range('D14').Value = '6.16'
$ws.Range('E14').Value = '  +0.52%  '
$ws.Range('D15').Value = '93.683.92'
$ws.Range('E15').Value = '  -4.01%  '
$ws.Range('D16').Value = '4.041.72'
$ws.Range('E16').Value = '  -0.52%  '
$ws.Range('D17').Value = '0.0000248'
$ws.Range('E17').Value = '  -3.90%  '
$ws.Range('D18').Value = '8.24'
$ws.Range('E18').Value = '  -10.32%  '
$ws.Range('D19').Value = '3.405.56'
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('D20').Value = '17.31'
$ws.Range('E20').Value = '  -5.13%  '
$ws.Range('D21').Value = '11.49'
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('D22').Value = '496.52'
$ws.Range('E22').Value = '  -3.10%  '
$ws.Range('D23').Value = '0.464'
$ws.Range('E23').Value = '  -11.16%  '
$ws.Range('E24').Value = '  -5.97%  '
$ws.Range('D25').Value = '0.0000190'
$ws.Range('E25').Value = '  -5.88%  '
$ws.Range('D26').Value = '6.44'
$ws.Range('E26').Value = '  -6.82%  '
$ws.Range('D27').Value = '90.77'
$ws.Range('E27').Value = '  -6.64%  '
$ws.Range('D28').Value = '3.590.33'
$ws.Range('E28').Value = '  +0.39%  '
$ws.Range('D29').Value = '11.79'
$ws.Range('E29').Value = '  -5.49%  '
$ws.Range('D30').Value = '11.40'
$ws.Range('E30').Value = '  -3.22%  '
$ws.Range('E31').Value = '  +0.19%  '
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('D33').Value = '0.135'
$ws.Range('E33').Value = '  -5.80%  '
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  +0.18%  '
$ws.Range('D35').Value = '0.174'
$ws.Range('E35').Value = '  -8.55%  '
$ws.Range('D36').Value = '29.26'
$ws.Range('E36').Value = '  +0.23%  '
$ws.Range('D37').Value = '0.543'
$ws.Range('E37').Value = '  -3.80%  '
$ws.Range('D38').Value = '535.66'
$ws.Range('E38').Value = '  +1.79%  '
$ws.Range('D39').Value = '7.62'
$ws.Range('E39').Value = '  -4.24%  '
$ws.Range('D40').Value = '1.43'
$ws.Range('E40').Value = '  -4.43%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').Value = '0.149'
$ws.Range('E42').Value = '  -2.64%  '
$ws.Range('D43').Value = '0.897'
$ws.Range('E43').Value = '  +3.79%  '
$ws.Range('E44').Value = '  -1.53%  '
$ws.Range('D45').Value = '3.70'
$ws.Range('E45').Value = '  +0.37%  '
$ws.Range('D46').Value = '1.69'
$ws.Range('E46').Value = '  -2.84%  '
$ws.Range('D47').Value = '5.58'
$ws.Range('E47').Value = '  -0.94%  '
$ws.Range('D48').Value = '2.18'
$ws.Range('E48').Value = '  -2.67%  '
$ws.Range('D49').Value = '3.26'
$ws.Range('E49').Value = '  +0.72%  '
$ws.Range('B50').Value = 'OKB'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D50').Value = '54.23'
$ws.Range('E50').Value = '  -2.88%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').Value = '0.0401'
$ws.Range('E51').Value = '  -5.87%  '
